# SwaadSutra_Consolidated_2026-01-13.xlsx update
# A new order (Sagar Borse / Til Poli x1) is inserted as the most recent
# order, so it becomes the new row 2 in "All Orders" and every existing
# order row shifts down by one. "Daily Summary" totals are bumped to
# reflect the extra order.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: "All Orders" ----
$ws = $wb.Worksheets.Item(1)

# Insert a brand-new blank row right under the header; existing rows 2-5
# (order #4, #3, #2, #1) shift down to rows 3-6 untouched.
$ws.Rows.Item(2).Insert()

# A couple of the new row's values look numeric/date-like to Excel's
# input parser ("7588930329", "2026-01-14") even though this column is
# text everywhere else in the sheet, so force those two cells to Text
# format before assigning them, keeping them as literal strings.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("J2").NumberFormat = "@"

$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "2026-01-13 16:40"
$ws.Range("C2").Value = "Sagar Borse"
$ws.Range("D2").Value = "A-1608"
$ws.Range("E2").Value = "7588930329"
$ws.Range("F2").Value = "Til Poli x1"
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = "NEW"
$ws.Range("I2").Value = "PENDING"
$ws.Range("J2").Value = "2026-01-14"
$ws.Range("K2").Value = "10:00"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""

# ---- Sheet 2: "Daily Summary" ----
$ws2 = $wb.Worksheets.Item(2)

# One more order today: Total Orders +1, Revenue/Pending +30 (the new
# order's total, still unpaid).
$ws2.Range("B2").Value = 5
$ws2.Range("E2").Value = 165
$ws2.Range("G2").Value = 165
